$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.3399353333333333
$ws.Range("H2").Value = 1.019806
$ws.Range("I2").Value = 0.09929991924017606
$ws.Range("J2").Value = 0.09929991924017606
$ws.Range("M2").Value = 7.487621999999999
$ws.Range("N2").Value = 22.462866
$ws.Range("O2").Value = 0.1384395179233961
$ws.Range("P2").Value = 0.1384395179233961
$ws.Range("Q2").Value = 2.545307280443999
$ws.Range("R2").Value = 22.907765523996
$ws.Range("S2").Value = 0.01374703294944214
$ws.Range("T2").Value = 0.01374703294944214
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.3399353333333333
$ws.Range("H3").Value = 1.019806
$ws.Range("I3").Value = 0.09929991924017606
$ws.Range("J3").Value = 0.09929991924017606
$ws.Range("O3").Value = 0.5916411627275552
$ws.Range("P3").Value = 0.5916411627275552
$ws.Range("Q3").Value = 10.877736223656
$ws.Range("R3").Value = 97.89962601290399
$ws.Range("S3").Value = 0.05874991967801009
$ws.Range("T3").Value = 0.05874991967801009
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.3399353333333333
$ws.Range("H4").Value = 1.019806
$ws.Range("I4").Value = 0.09929991924017606
$ws.Range("J4").Value = 0.09929991924017606
$ws.Range("M4").Value = 14.59882166666667
$ws.Range("N4").Value = 43.796465
$ws.Range("O4").Value = 0.2699193193490487
$ws.Range("P4").Value = 0.2699193193490487
$ws.Range("Q4").Value = 4.962655309532222
$ws.Range("R4").Value = 44.66389778579
$ws.Range("S4").Value = 0.02680296661272383
$ws.Range("T4").Value = 0.02680296661272383
$ws.Range("I5").Value = 0.4094685684206303
$ws.Range("J5").Value = 0.4094685684206303
$ws.Range("M5").Value = 7.487621999999999
$ws.Range("N5").Value = 22.462866
$ws.Range("O5").Value = 0.1384395179233961
$ws.Range("P5").Value = 0.1384395179233961
$ws.Range("Q5").Value = 10.49571174165
$ws.Range("R5").Value = 94.46140567484998
$ws.Range("S5").Value = 0.05668663121693517
$ws.Range("T5").Value = 0.05668663121693518
$ws.Range("I6").Value = 0.4094685684206303
$ws.Range("J6").Value = 0.4094685684206303
$ws.Range("O6").Value = 0.5916411627275552
$ws.Range("P6").Value = 0.5916411627275552
$ws.Range("S6").Value = 0.2422584599207692
$ws.Range("T6").Value = 0.2422584599207692
$ws.Range("I7").Value = 0.4094685684206303
$ws.Range("J7").Value = 0.4094685684206303
$ws.Range("M7").Value = 14.59882166666667
$ws.Range("N7").Value = 43.796465
$ws.Range("O7").Value = 0.2699193193490487
$ws.Range("P7").Value = 0.2699193193490487
$ws.Range("Q7").Value = 20.46377661440278
$ws.Range("R7").Value = 184.173989529625
$ws.Range("S7").Value = 0.1105234772829259
$ws.Range("T7").Value = 0.1105234772829259
$ws.Range("G8").Value = 1.681642333333333
$ws.Range("H8").Value = 5.044927
$ws.Range("I8").Value = 0.4912315123391937
$ws.Range("J8").Value = 0.4912315123391937
$ws.Range("M8").Value = 7.487621999999999
$ws.Range("N8").Value = 22.462866
$ws.Range("O8").Value = 0.1384395179233961
$ws.Range("P8").Value = 0.1384395179233961
$ws.Range("Q8").Value = 12.591502131198
$ws.Range("R8").Value = 113.323519180782
$ws.Range("S8").Value = 0.06800585375701876
$ws.Range("T8").Value = 0.06800585375701876
$ws.Range("G9").Value = 1.681642333333333
$ws.Range("H9").Value = 5.044927
$ws.Range("I9").Value = 0.4912315123391937
$ws.Range("J9").Value = 0.4912315123391937
$ws.Range("O9").Value = 0.5916411627275552
$ws.Range("P9").Value = 0.5916411627275552
$ws.Range("Q9").Value = 53.811592767252
$ws.Range("R9").Value = 484.304334905268
$ws.Range("S9").Value = 0.2906327831287759
$ws.Range("T9").Value = 0.2906327831287759
$ws.Range("G10").Value = 1.681642333333333
$ws.Range("H10").Value = 5.044927
$ws.Range("I10").Value = 0.4912315123391937
$ws.Range("J10").Value = 0.4912315123391937
$ws.Range("M10").Value = 14.59882166666667
$ws.Range("N10").Value = 43.796465
$ws.Range("O10").Value = 0.2699193193490487
$ws.Range("P10").Value = 0.2699193193490487
$ws.Range("Q10").Value = 24.54999653145056
$ws.Range("R10").Value = 220.949968783055
$ws.Range("S10").Value = 0.132592875453399
$ws.Range("T10").Value = 0.132592875453399
